$d = $word.ActiveDocument
$tab = "`t"

# --- Hunk 1 -----------------------------------------------------------
# Merge the two runs holding "SUN Apr 29" and " 11:11:14 IST 2018" into
# a single run with the combined text (mirrors the diff's run merge).
$mergeRange = $d.Content
$mergeRange.Find.Execute("SUN Apr 29 11:11:14 IST 2018", $false, $false, `
  $false, $false, $false, $true, 1, $false, `
  "SUN Apr 29 11:11:14 IST 2018", 2) | Out-Null

# --- Hunk 2 -----------------------------------------------------------
# Append a new "MON APR 30" purchase record right after the existing
# "Amount balance ... - 353.0" paragraph (and before the blank
# paragraphs that close out the document).
$paras = $d.Paragraphs
$targetIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
  if ($paras.Item($i).Range.Text -like "*- 353.0*") {
    $targetIndex = $i
  }
}

$insertionPoint = $paras.Item($targetIndex).Range
$insertionPoint.Collapse(0)

$newBlock = "`r" + `
  "MON APR 30 10:39:53 IST 2018`r" + `
  "Person Name" + $tab + $tab + $tab + $tab + "- EERANNA`r" + `
  "---------------------------------------------------------------`r" + `
  "Item Name" + $tab + $tab + $tab + $tab + "- POTATO`r" + `
  "Amount Received" + $tab + $tab + $tab + "- 200`r" + `
  "Amount balance" + $tab + $tab + $tab + "- 153.0`r" + `
  "Amount Received mode" + $tab + $tab + "- CASH`r" + `
  "`r" + `
  "`r"

$insertionPoint.InsertAfter($newBlock)

# Re-fetch the paragraph collection (it grew) and fix up per-paragraph
# formatting so it matches the source record's look:
#   targetIndex+1  -> blank line                         (bold, like the record separator)
#   targetIndex+2  -> "MON APR 30 ..." timestamp line     (plain)
#   targetIndex+3  -> "Person Name ... - EERANNA"         (plain)
#   targetIndex+4  -> "----...----" divider               (plain)
#   targetIndex+5  -> "Item Name ... - POTATO"            (plain)
#   targetIndex+6  -> "Amount Received ... - 200"         (plain, red)
#   targetIndex+7  -> "Amount balance ... - 153.0"        (bold)
#   targetIndex+8  -> "Amount Received mode ... - CASH"   (plain)
#   targetIndex+9  -> blank line                          (plain)
#   targetIndex+10 -> blank line                          (bold)
$paras = $d.Paragraphs

$pBlank1 = $paras.Item($targetIndex + 1)
$pBlank1.Range.Font.Bold = 1

$pDate = $paras.Item($targetIndex + 2)
$pDate.Range.Font.Bold = 0

$pPerson = $paras.Item($targetIndex + 3)
$pPerson.Range.Font.Bold = 0

$pDivider = $paras.Item($targetIndex + 4)
$pDivider.Range.Font.Bold = 0

$pItem = $paras.Item($targetIndex + 5)
$pItem.Range.Font.Bold = 0

$pReceived = $paras.Item($targetIndex + 6)
$pReceived.Range.Font.Bold = 0
$pReceived.Range.Font.Color = 255

$pBalance = $paras.Item($targetIndex + 7)
$pBalance.Range.Font.Bold = 1

$pMode = $paras.Item($targetIndex + 8)
$pMode.Range.Font.Bold = 0

$pBlank2 = $paras.Item($targetIndex + 9)
$pBlank2.Range.Font.Bold = 0

$pBlank3 = $paras.Item($targetIndex + 10)
$pBlank3.Range.Font.Bold = 1

Write-Host "Done. Paragraphs: $($d.Paragraphs.Count)"
